$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.806427121162415
$ws.Range("B1").Value = 3.581918239593506
$ws.Range("C1").Value = 3.147284269332886
$ws.Range("D1").Value = 3.414026737213135
$ws.Range("E1").Value = 1.729030251502991
